# suat-an.xlsx update: add "THEO DOI SUAT AN" title row, format header block
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new title row above the existing header rows ---
$null = $ws.Rows(1).Insert()

# --- Title cell (merged A1:N1) ---
$titleRange = $ws.Range("A1:N1")
$null = $titleRange.Merge()
$ws.Range("A1").Value = "THEO DÕI SUẤT ĂN"

# --- Column widths (approximate Excel's best-fit widths for the final layout) ---
$ws.Columns("A").ColumnWidth = 22.0
$ws.Columns("B").ColumnWidth = 22.666666666666668
$ws.Columns("C").ColumnWidth = 14.5
$ws.Columns("D").ColumnWidth = 11.666666666666666
$ws.Columns("E").ColumnWidth = 14.5
$ws.Columns("F").ColumnWidth = 11.666666666666666
$ws.Columns("G").ColumnWidth = 14.5
$ws.Columns("H").ColumnWidth = 11.666666666666666
$ws.Columns("I").ColumnWidth = 14.5
$ws.Columns("J").ColumnWidth = 11.666666666666666
$ws.Columns("K").ColumnWidth = 14.5
$ws.Columns("L").ColumnWidth = 11.666666666666666
$ws.Columns("M").ColumnWidth = 14.5
$ws.Columns("N").ColumnWidth = 11.666666666666666

# --- Row heights for the three header rows ---
$ws.Rows("1:3").RowHeight = 22.5

# --- Fonts, fill, borders, alignment for the whole header block (rows 1-3) ---
$headerBlock = $ws.Range("A1:N3")
$headerBlock.Font.Name = "Times New Roman"
$headerBlock.Font.Family = 1
$headerBlock.Interior.Color = 14277081
$headerBlock.Borders.LineStyle = 1
$headerBlock.Borders.Weight = 2
$headerBlock.HorizontalAlignment = -4108

$row23 = $ws.Range("A2:N3")
$row23.VerticalAlignment = -4108

$titleRange.Font.Size = 18
$titleRange.Font.Bold = $true

$row23Font = $ws.Range("A2:N3")
$row23Font.Font.Size = 18
$row23Font.Font.Bold = $true

$dataFont = $ws.Range("A1:N1")
$dataFont.Font.Size = 18

# Row 2/3 (original header rows) retain the same big bold font as the title per source file
# --- Selection matches the saved cursor position in the source workbook ---
$null = $ws.Range("B6").Select()
